$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Ababei Vasile
$ws.Range("B2").Value = 2
$ws.Range("D2").Value = "fibonacci"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = "Input: 4 Expected Output: 3 Actual Output: 3Input: 5 Expected Output: 5 Actual Output: 5Success!"
$ws.Range("H2").Value = 0.5
$ws.Range("I2").Value = " The function respects the requirement because it has a single loop that iterates 'n' times, resulting in a time complexity of O(n)."
$ws.Range("J2").Value = 0.5
$ws.Range("K2").Value = " The function respects the requirement because it only uses a constant amount of space to store the variables fib, prevFib, and temp, which does not grow with the input size n, hence the space complexity is O(1) which is at most O(n)."
$ws.Range("O2").Value = " The function does not handle incorrect input (negative numbers) as it does not check for negative values of n and will enter an infinite loop if n is negative."

# Row 3 - Alexe Robert George
$ws.Range("B3").Value = 1
$ws.Range("D3").Value = "fibonacci"
$ws.Range("E3").Value = 1
$ws.Range("G3").Value = "Input: 4 Expected Output: 3 Actual Output: 1"
$ws.Range("I3").Value = " The function's time complexity is not O(n) because the line ""fib *= prevFib;"" inside the loop causes the numbers to grow exponentially, leading to an exponential time complexity."
$ws.Range("J3").Value = 0.5
$ws.Range("K3").Value = " The function respects the requirement because it only uses a constant amount of space to store the variables fib, prevFib, and temp, which does not grow with the input size n, hence the space complexity is O(1) which is at most O(n)."
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = "Input: 0 Expected Output: 0 Actual Output: 0Input: 1 Expected Output: 1 Actual Output: 1Success!"
$ws.Range("O3").Value = " The function does not handle incorrect input (negative numbers and non-integers) as it does not include any error checking or handling for such cases."

# Row 4 - Alin Claudiu
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "Yes"
$ws.Range("D4").Value = "fibonacci"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = "Input: 4 Expected Output: 3 Actual Output: 3Input: 5 Expected Output: 5 Actual Output: 5Success!"
$ws.Range("H4").Value = 0.5
$ws.Range("I4").Value = " The function respects the requirement because it has a single loop that iterates 'n' times, resulting in a time complexity of O(n)."
$ws.Range("J4").Value = 0.5
$ws.Range("K4").Value = " The function respects the requirement because it only uses a constant amount of space to store the variables fib, prevFib, and temp, which does not grow with the input size n, hence the space complexity is O(1) which is at most O(n)."
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = "Input: 0 Expected Output: 0 Actual Output: 0Input: 1 Expected Output: 1 Actual Output: 1Success!"
$ws.Range("O4").Value = " The function does not handle incorrect input (negative numbers and non-integers) as it does not include any error checking or handling for such cases."

# Column widths grew because of the new (longer) text content in columns G, I, K and O.
$ws.Columns.Item(7).ColumnWidth = 85.83333333333334
$ws.Columns.Item(9).ColumnWidth = 156.16666666666666
$ws.Columns.Item(11).ColumnWidth = 204.66666666666666
$ws.Columns.Item(15).ColumnWidth = 135.83333333333331

Write-Host "Done"
